$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.700.94'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '1.893.24'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -1.13%  '
$ws.Range("D5").Value = '''311.83'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").Value = '''0.4881'
$ws.Range("E7").Value = '  +1.08%  '
$ws.Range("D8").Value = '''0.3791'
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("D9").Value = '''0.07325'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").Value = '''0.9121'
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("D11").Value = '''20.56'
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").Value = '''0.07648'
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("D13").Value = '1.902.61'
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").Value = '''5.474'
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").Value = '''6.619'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '''91.29'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("D18").Value = '''0.000008764'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").Value = '27.791.13'
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").Value = '''14.47'
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").Value = '''5.117'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '2.140.00'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''154.11'
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '''1.899'
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("D27").Value = '''18.36'
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("D28").Value = '''2.156'
$ws.Range("E28").Value = '  +5.71%  '
$ws.Range("D29").Value = '''115.30'
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("D30").Value = '''4.881'
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").Value = '''0.08904'
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").Value = '''0.7677'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("D35").Value = '''4.630'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").Value = '''2.565'
$ws.Range("E36").Value = '  -6.05%  '
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("D38").Value = '''1.095'
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("D40").Value = '''0.5469'
$ws.Range("E40").Value = '  -2.63%  '
$ws.Range("D41").Value = '''2.981'
$ws.Range("E41").Value = '  -0.67%  '
$ws.Range("D42").Value = '''6.897'
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("D43").Value = '''8.507'
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''112.53'
$ws.Range("E44").Value = '  +6.32%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '''0.1518'
$ws.Range("E45").Value = '  -0.68%  '
$ws.Range("D46").Value = '''10.65'
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '''0.4786'
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D49").Value = '''1.640'
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("D50").Value = '''67.35'
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("D51").Value = '''0.06052'
$ws.Range("E51").Value = '  -1.11%  '
